{"js": "// 1. Update the \"Curso (semestre ideal)\" line: drop the \"EB (8), \" part.\nconst hits = context.document.body.search(\"EB (8), EQD (10), EQN (12)\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"EQD (10), EQN (12)\", \"Replace\");\n} else {\n  // Fallback: narrower search in case formatting/spacing differs slightly.\n  const hits2 = context.document.body.search(\"EB (8), \", { matchCase: true });\n  hits2.load(\"text\");\n  await context.sync();\n  if (hits2.items.length > 0) {\n    hits2.items[0].insertText(\"\", \"Replace\");\n  }\n}\nawait context.sync();\n\n// 2. Remove the \"Requisitos\" heading paragraph and the bullet paragraph\n//    right after it (the \"LOQ4046 - ...\" requirement line).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\" && items[i].text.trim() === \"Requisitos\") {\n    // Delete the following paragraph first (if present) so indices of\n    // earlier paragraphs stay valid, then delete the heading itself.\n    if (i + 1 < items.length) {\n      items[i + 1].delete();\n    }\n    items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"Curso (semestre ideal): EB (8), EQD (10), EQN (12)\" ->\n#    \"Curso (semestre ideal): EQD (10), EQN (12)\"\n#    Drop the \"EB (8), \" portion via Find & Replace on the whole document.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"EB (8), \"\n$find.Replacement.Text = \"\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 2. Remove the \"Requisitos\" heading and the requirement bullet right after\n#    it (the \"LOQ4046 - ...\" line), which together made up the whole\n#    \"Requisitos\" section at the end of the document.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.Trim() -eq \"Requisitos\") {\n        if ($i -lt $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($i + 1).Range.Delete()\n        }\n        $p.Range.Delete()\n        break\n    }\n}\n"}
